$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.210.32'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range('E2').Value = '  -1.06%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.659.12'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range('E3').Value = '  -0.74%  '

$ws.Range('E4').Value = '  +0.26%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.29'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range('E5').Value = '  -1.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5165'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range('E6').Value = '  -2.20%  '

$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2641'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range('E8').Value = '  -1.38%  '

$ws.Range('E9').Value = '  -1.63%  '

$ws.Range('E10').Value = '  -4.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07761'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range('E11').Value = '  -0.61%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.484'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range('E12').Value = '  -0.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.629.17'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range('E13').Value = '  -2.39%  '

$ws.Range('E14').Value = '  -0.76%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5456'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range('E15').Value = '  -2.14%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8130'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Range('E16').Value = '  -2.07%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.81'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range('E17').Value = '  -1.12%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.215.14'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('E19').Value = '  +0.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.606'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range('E20').Value = '  -3.29%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.00'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range('E21').Value = '  -0.63%  '

$ws.Range('E22').Value = '  -2.25%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.989'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range('E23').Value = '  -5.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.006'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range('E24').Value = '  +0.36%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '139.44'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range('E25').Value = '  +0.44%  '

$ws.Range('E26').Value = '  -3.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.273'
$ws.Cells.Item(27, 4).Style = "Normal"

$ws.Range('E28').Value = '  -1.16%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.443'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range('E29').Value = '  +1.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05929'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range('E30').Value = '  -4.49%  '

$ws.Range('E31').Value = '  -1.12%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.547'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range('E32').Value = '  -1.71%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.276'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Range('E33').Value = '  -4.29%  '

$ws.Range('E34').Value = '  -6.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9588'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range('E35').Value = '  -4.78%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.421'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range('E36').Value = '  +0.25%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.770'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range('E37').Value = '  -0.34%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5662'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range('E38').Value = '  -7.69%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.042'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range('E39').Value = '  +0.11%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01589'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range('E40').Value = '  -1.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8556'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range('E41').Value = '  +0.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.003'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range('E42').Value = '  +0.24%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.012.39'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range('E43').Value = '  -7.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.63'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range('E44').Value = '  +0.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.800.43'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range('E45').Value = '  -0.86%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈108'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range('E46').Value = '  -3.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.46'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range('E47').Value = '  -3.08%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Range('E48').Value = '  +0.38%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.051'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range('E49').Value = '  -0.57%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05166'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range('E50').Value = '  -0.54%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4216'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Range('E51').Value = '  -0.31%  '

